$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update metric values
$ws.Range("B3").Value = 0.9994150917987699
$ws.Range("C3").Value = 0.9993859305186326
$ws.Range("D3").Value = 0.9983299113817082

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9994282868642346
$ws.Range("C4").Value = 0.9994351815338799
$ws.Range("D4").Value = 0.9981299281825177

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9994905602732072
$ws.Range("C5").Value = 0.9994696434829643
$ws.Range("D5").Value = 0.9995420842366699
